# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-11-19 Wednesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-11-20 Thursday", 2)

# Update each division-problem cell in the single table by explicit
# (row, column) address, so that cells whose new value happens to equal
# another cell's old value can never collide during the edit.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "87÷7=" },
    @{ Row = 1;  Col = 2; Text = "30÷5=" },
    @{ Row = 1;  Col = 3; Text = "87÷9=" },
    @{ Row = 1;  Col = 4; Text = "18÷2=" },
    @{ Row = 1;  Col = 5; Text = "17÷4=" },

    @{ Row = 5;  Col = 1; Text = "51÷9=" },
    @{ Row = 5;  Col = 2; Text = "76÷5=" },
    @{ Row = 5;  Col = 3; Text = "52÷7=" },
    @{ Row = 5;  Col = 4; Text = "25÷5=" },
    @{ Row = 5;  Col = 5; Text = "21÷8=" },

    @{ Row = 9;  Col = 1; Text = "68÷3=" },
    @{ Row = 9;  Col = 2; Text = "86÷4=" },
    @{ Row = 9;  Col = 3; Text = "60÷9=" },
    @{ Row = 9;  Col = 4; Text = "72÷6=" },
    @{ Row = 9;  Col = 5; Text = "33÷2=" },

    @{ Row = 13; Col = 1; Text = "94÷7=" },
    @{ Row = 13; Col = 2; Text = "32÷3=" },
    @{ Row = 13; Col = 3; Text = "31÷6=" },
    @{ Row = 13; Col = 4; Text = "74÷8=" },
    @{ Row = 13; Col = 5; Text = "27÷8=" },

    @{ Row = 17; Col = 1; Text = "19÷7=" },
    @{ Row = 17; Col = 2; Text = "60÷2=" },
    @{ Row = 17; Col = 3; Text = "25÷2=" },
    @{ Row = 17; Col = 4; Text = "14÷9=" },
    @{ Row = 17; Col = 5; Text = "13÷2=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
